$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Available Data column (B2:B5) from 313 to 316
$ws.Range("B2").Value = 316
$ws.Range("B3").Value = 316
$ws.Range("B4").Value = 316
$ws.Range("B5").Value = 316

# Update rows 6-8 (Oporavljeni, Testirani, Smrtni sl.)
$ws.Range("B6").Value = 238
$ws.Range("C6").Value = 78
$ws.Range("D6").Value = 0.3277310924369748

$ws.Range("B7").Value = 238
$ws.Range("C7").Value = 78
$ws.Range("D7").Value = 0.3277310924369748

$ws.Range("B8").Value = 238
$ws.Range("C8").Value = 78
$ws.Range("D8").Value = 0.3277310924369748
